$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = "INDONESIA"
$ws.Range("B36").Value = 7717
$ws.Range("C36").Value = 6759
$ws.Range("D36").Value = 5972
$ws.Range("E36").Value = 5705
$ws.Range("F36").Value = 5158

$ws.Range("B1").Select()
